$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that needs to move from
# 45188 (2023-09-19) to 45189 (2023-09-20) for every data row (2-115).
for ($r = 2; $r -le 115; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
